$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("792÷6=", $true, $false, $false, $false, $false, $true, 0, $false, "651÷6=", 2)

$cell = $t.Cell(1, 2)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("744÷6=", $true, $false, $false, $false, $false, $true, 0, $false, "553÷7=", 2)

$cell = $t.Cell(1, 3)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("264÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "973÷7=", 2)

$cell = $t.Cell(1, 4)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("744÷6=", $true, $false, $false, $false, $false, $true, 0, $false, "206÷6=", 2)

$cell = $t.Cell(1, 5)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("260÷9=", $true, $false, $false, $false, $false, $true, 0, $false, "219÷2=", 2)

$cell = $t.Cell(5, 1)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("789÷4=", $true, $false, $false, $false, $false, $true, 0, $false, "469÷9=", 2)

$cell = $t.Cell(5, 2)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("803÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "905÷2=", 2)

$cell = $t.Cell(5, 3)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("518÷4=", $true, $false, $false, $false, $false, $true, 0, $false, "360÷4=", 2)

$cell = $t.Cell(5, 4)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("379÷3=", $true, $false, $false, $false, $false, $true, 0, $false, "793÷5=", 2)

$cell = $t.Cell(5, 5)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("446÷5=", $true, $false, $false, $false, $false, $true, 0, $false, "532÷8=", 2)

$cell = $t.Cell(9, 1)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("477÷3=", $true, $false, $false, $false, $false, $true, 0, $false, "625÷9=", 2)

$cell = $t.Cell(9, 2)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("537÷2=", $true, $false, $false, $false, $false, $true, 0, $false, "962÷7=", 2)

$cell = $t.Cell(9, 3)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("787÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "751÷8=", 2)

$cell = $t.Cell(9, 4)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("561÷4=", $true, $false, $false, $false, $false, $true, 0, $false, "863÷2=", 2)

$cell = $t.Cell(9, 5)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("103÷2=", $true, $false, $false, $false, $false, $true, 0, $false, "682÷3=", 2)

$cell = $t.Cell(13, 1)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("662÷9=", $true, $false, $false, $false, $false, $true, 0, $false, "474÷7=", 2)

$cell = $t.Cell(13, 2)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("937÷3=", $true, $false, $false, $false, $false, $true, 0, $false, "251÷4=", 2)

$cell = $t.Cell(13, 3)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("274÷6=", $true, $false, $false, $false, $false, $true, 0, $false, "908÷2=", 2)

$cell = $t.Cell(13, 4)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("723÷2=", $true, $false, $false, $false, $false, $true, 0, $false, "647÷9=", 2)

$cell = $t.Cell(13, 5)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("560÷2=", $true, $false, $false, $false, $false, $true, 0, $false, "988÷2=", 2)

$cell = $t.Cell(17, 1)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("514÷7=", $true, $false, $false, $false, $false, $true, 0, $false, "170÷3=", 2)

$cell = $t.Cell(17, 2)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("956÷7=", $true, $false, $false, $false, $false, $true, 0, $false, "924÷4=", 2)

$cell = $t.Cell(17, 3)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("592÷4=", $true, $false, $false, $false, $false, $true, 0, $false, "642÷6=", 2)

$cell = $t.Cell(17, 4)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("681÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "574÷9=", 2)

$cell = $t.Cell(17, 5)
$rng = $d.Range($cell.Range.Start, $cell.Range.End)
$rng.Find.Execute("837÷7=", $true, $false, $false, $false, $false, $true, 0, $false, "527÷7=", 2)
